$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "Notes" (column G) entries describing measurement issues found 03.04.24 ---

# Row 12: "Distance Average Affinity" / affinity / Affinity - append info to the existing note
$ws.Range("G12").Value = "03.04.24: `nMeasure cant be calculated with the provided tool, throws an exception`nUnclear if  ""Affinity"" is the same as ""Distance Average Affinity"""

# Row 11: "Number and Percentage of Unique Sequences" / distinct_traces
$ws.Range("G11").Value = "03.04.24: `nMeasure cant be calculated with the provided tool, throws an exception`n"

# Row 20: "Enriched Variant entropy"
$ws.Range("G20").Value = "03.04.24:`nHow to calculate the enriched variant entropy? With what is the normal variant entropy enriched?"

# Row 24: "Number of Activity Repetitions"
$ws.Range("G24").Value = "03.04.24:`nThe authors state that this measure can be calculated out of the box by PM4PY. However, I only found the method ""pm4py.get_rework_cases_per_activity(log)"" and this method only returns a dict of reps per process step. How to turn this into a valid metric?"

# Row 30: "structure"
$ws.Range("G30").Value = "03.04.24: `nMeasure cant be calculated with the provided tool, throws an exception`n"

# Row 31: "Intercept"
$ws.Range("G31").Value = "03.04.24:`nWhat is this measure and how can it be computed?"

# --- Highlight the metrics affected by the new issues in red (matches other problematic metrics) ---
$ws.Range("E11").Interior.Color = 255
$ws.Range("E12").Interior.Color = 255
$ws.Range("E30").Interior.Color = 255

# --- Keep the saved view roughly where the author left off scrolling/selecting ---
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("I23").Select()
